{"js": "// Replace the accented placeholder \"{{direcci\u00f3n}}\" with the accent-free\n// \"{{direccion}}\" in the body of the document (forms template clean-up).\nconst body = context.document.body;\n\nconst results = body.search(\"{{direcci\u00f3n}}\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"{{direccion}}\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the accented placeholder \"{{direcci\u00f3n}}\" with the accent-free\n# \"{{direccion}}\" in the document body (forms template clean-up).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Execute(\n    \"{{direcci\u00f3n}}\",   # FindText\n    $true,             # MatchCase\n    $false,            # MatchWholeWord\n    $false,            # MatchWildcards\n    $false,            # MatchSoundsLike\n    $false,            # MatchAllWordForms\n    $true,             # Forward\n    1,                 # Wrap (wdFindContinue)\n    $false,            # Format\n    \"{{direccion}}\",   # ReplaceWith\n    2                  # Replace (wdReplaceAll)\n)\n"}
